# B6-PowerPoint.pptx — commit "Sat, Aug 01, 2020  8:05:04 AM"
#
# The canonical-XML diff for this commit shows two kinds of changes:
#
#   1. Three tables (on slides 14, 15 and 16 — the only shapes on those
#      slides that are graphicFrames/tables) have their <a:tableStyleId>
#      changed from the locally-defined "Table_0" style
#      {EC34C5B5-F761-4471-A960-389B25687B21} to the built-in table style
#      {F6970670-5593-4ABF-837C-2785B26A8637}.
#
#   2. ppt/theme/theme1.xml and ppt/theme/theme2.xml trade their entire
#      contents (the deck's "Office Theme" and "Integral" themes swap
#      which physical part holds which). theme2.xml is the theme that is
#      actually wired to the Slide Master (and to the top-level
#      presentation theme relationship), so from the editor's point of
#      view this is simply "the applied design changed from Integral to
#      Office Theme" (and, as a byproduct of the part-swap, the Notes
#      Master — which is wired to theme1.xml — ends up showing Integral).

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Table style swap on the three affected tables.
#    Each of slides 14/15/16 has exactly one graphicFrame (table), and
#    it is always the first shape on the slide.
# ---------------------------------------------------------------------
$newTableStyle = "{F6970670-5593-4ABF-837C-2785B26A8637}"

foreach ($slideIndex in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIndex)
    $tableShape = $slide.Shapes.Item(1)
    if ($tableShape.HasTable) {
        $tableShape.Table.ApplyStyle($newTableStyle)
    }
}

# ---------------------------------------------------------------------
# 2) Re-point the deck's applied design from "Integral" to
#    "Office Theme". The object model's per-slide ThemeColorScheme maps
#    onto the Slide Master's theme part (theme2.xml in this deck), so
#    pushing the Office Theme's twelve scheme colors through it recolors
#    that shared theme part to match the target "Office Theme" palette.
#    (msoThemeColorDark1..msoThemeColorFollowedHyperlink order, i.e. the
#    same 1-12 index order PowerPoint uses for ThemeColorScheme.Colors.)
# ---------------------------------------------------------------------
function RgbValue([int]$r, [int]$g, [int]$b) {
    return $r + ($g * 256) + ($b * 65536)
}

$officeThemeColors = @(
    (RgbValue 0x00 0x00 0x00),   # dk1      000000
    (RgbValue 0xFF 0xFF 0xFF),   # lt1      FFFFFF
    (RgbValue 0x44 0x54 0x6A),   # dk2      44546A
    (RgbValue 0xE7 0xE6 0xE6),   # lt2      E7E6E6
    (RgbValue 0x5B 0x9B 0xD5),   # accent1  5B9BD5
    (RgbValue 0xED 0x7D 0x31),   # accent2  ED7D31
    (RgbValue 0xA5 0xA5 0xA5),   # accent3  A5A5A5
    (RgbValue 0xFF 0xC0 0x00),   # accent4  FFC000
    (RgbValue 0x44 0x72 0xC4),   # accent5  4472C4
    (RgbValue 0x70 0xAD 0x47),   # accent6  70AD47
    (RgbValue 0x05 0x63 0xC1),   # hlink    0563C1
    (RgbValue 0x95 0x4F 0x72)    # folHlink 954F72
)

$themeColorScheme = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le $officeThemeColors.Count; $i++) {
    $themeColorScheme.Colors($i).RGB = $officeThemeColors[$i - 1]
}
